$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F5").Value = 4
$ws.Range("F12").Value = 1
$ws.Range("F18").Value = 2
$ws.Range("F23").Value = 2
$ws.Range("F24").Value = -2
$ws.Range("E33").Value = 3
$ws.Range("F33").Value = -2
$ws.Range("F38").Value = -1
$ws.Range("F39").Value = -2
$ws.Range("F43").Value = 0
$ws.Range("F53").Value = 0
$ws.Range("F54").Value = -1
